$wb = $excel.ActiveWorkbook

# NOTE: sheets "Vector_bf" (index 5) and "Vector_BF" (index 6) differ only by
# case, and Worksheets.Item(name) lookup is case-insensitive, so we must use
# 1-based positional indices (matching workbook.xml sheet order) to address
# them unambiguously.
$ws1 = $wb.Worksheets.Item(1)   # Funciones_Objetivo
$ws2 = $wb.Worksheets.Item(2)   # Restricciones_del_lider
$ws3 = $wb.Worksheets.Item(3)   # Restricciones_del_follower
$ws4 = $wb.Worksheets.Item(4)   # Punto_modificado
$ws5 = $wb.Worksheets.Item(5)   # Vector_bf
$ws6 = $wb.Worksheets.Item(6)   # Vector_BF
$ws7 = $wb.Worksheets.Item(7)   # Vector_Alpha

function Set-TextValue($range, $value) {
    # Force Excel to store the value as text (shared string) rather than
    # auto-converting numeric-looking strings into numbers, while keeping
    # the cell style back to the default (no style index written).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Sheet2: Restricciones_del_lider
Set-TextValue $ws2.Range("A2") "0.049999999999998934 - x + y"
Set-TextValue $ws2.Range("B2") "-0.049999999999998934"
Set-TextValue $ws2.Range("D2") "0.4"

# Sheet3: Restricciones_del_follower
Set-TextValue $ws3.Range("A2") "-230.45 + x + 15.0y"
Set-TextValue $ws3.Range("B2") "210.45"
Set-TextValue $ws3.Range("D2") "0.55"
Set-TextValue $ws3.Range("E2") "7.800000000000001"
Set-TextValue $ws3.Range("F2") "9.9"

Set-TextValue $ws3.Range("A3") "4.896 - 0.33999999999999997y"
Set-TextValue $ws3.Range("B3") "-4.896"
Set-TextValue $ws3.Range("D3") "0.45"
Set-TextValue $ws3.Range("E3") "8.4"
Set-TextValue $ws3.Range("F3") "3.1"

Set-TextValue $ws3.Range("A4") "-103.31545454545453 + 7.12121212121212y"
Set-TextValue $ws3.Range("B4") "82.54545454545453"
Set-TextValue $ws3.Range("D4") "0.55"
Set-TextValue $ws3.Range("E4") "7.7"
Set-TextValue $ws3.Range("F4") "4.699999999999999"

# Sheet4: Punto_modificado
Set-TextValue $ws4.Range("A2") "14.45"
Set-TextValue $ws4.Range("B2") "14.4"

# Sheet5: Vector_bf
Set-TextValue $ws5.Range("A2") "-65.01366666666667"

# Sheet6: Vector_BF
Set-TextValue $ws6.Range("A2") "-38.94"
Set-TextValue $ws6.Range("A3") "-183.45733333333334"

# Sheet7: Vector_Alpha (A2 is a real numeric literal cell, not a shared string)
$ws7.Range("A2").Value = 0.66
